$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("K2").Value = 1618
$ws.Range("K3").Value = 1536
$ws.Range("F4").Value = 1907
$ws.Range("K4").Value = 336
$ws.Range("K6").Value = 2010
$ws.Range("F7").Value = 24100
$ws.Range("K7").Value = 5597

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("K2").Value = 41
$ws.Range("K6").Value = 46
$ws.Range("K7").Value = 158
$ws.Range("K8").Value = 352
$ws.Range("K11").Value = 117
$ws.Range("K19").Value = 151
$ws.Range("K20").Value = 128
$ws.Range("K22").Value = 13
$ws.Range("K25").Value = 25
$ws.Range("K29").Value = 255
$ws.Range("K30").Value = 20
$ws.Range("K33").Value = 230
$ws.Range("K36").Value = 63
$ws.Range("K37").Value = 190
$ws.Range("K41").Value = 59
$ws.Range("K42").Value = 203
$ws.Range("K43").Value = 58
$ws.Range("K44").Value = 50
$ws.Range("K45").Value = 6
$ws.Range("K47").Value = 39
$ws.Range("K52").Value = 147
$ws.Range("K53").Value = 88
$ws.Range("K55").Value = 58
$ws.Range("K60").Value = 43
$ws.Range("F63").Value = 193
$ws.Range("K63").Value = 18
$ws.Range("K64").Value = 35
$ws.Range("K65").Value = 142
$ws.Range("K67").Value = 219
$ws.Range("K68").Value = 15
$ws.Range("K71").Value = 16
$ws.Range("K75").Value = 20
$ws.Range("K76").Value = 83
$ws.Range("K79").Value = 151
$ws.Range("K82").Value = 9
$ws.Range("K83").Value = 113
$ws.Range("K85").Value = 284
$ws.Range("K91").Value = 53
$ws.Range("K92").Value = 25
$ws.Range("K94").Value = 68
$ws.Range("K97").Value = 47
$ws.Range("K98").Value = 37
$ws.Range("K99").Value = 100
$ws.Range("K100").Value = 8
$ws.Range("F101").Value = 24100
$ws.Range("K101").Value = 5597

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("K6").Value = 38
$ws.Range("K7").Value = 158

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("K2").Value = 40
$ws.Range("K7").Value = 117

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("K2").Value = 103
$ws.Range("K4").Value = 16
$ws.Range("K7").Value = 284

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("K3").Value = 30
$ws.Range("K7").Value = 147

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("K2").Value = 19
$ws.Range("K6").Value = 51
$ws.Range("K7").Value = 88

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("K3").Value = 103
$ws.Range("K6").Value = 120
$ws.Range("K7").Value = 352

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("K4").Value = 6
$ws.Range("K7").Value = 113

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("K3").Value = 85
$ws.Range("K7").Value = 230

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("K2").Value = 41
$ws.Range("K6").Value = 70
$ws.Range("K7").Value = 190

$ws = $wb.Worksheets.Item("New City")
$ws.Range("K2").Value = 40
$ws.Range("K6").Value = 62
$ws.Range("K7").Value = 142

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("K2").Value = 33
$ws.Range("K6").Value = 29
$ws.Range("K7").Value = 100

$ws = $wb.Worksheets.Item("Fuller Park")
$ws.Range("K2").Value = 4
$ws.Range("K7").Value = 20

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("K3").Value = 70
$ws.Range("K6").Value = 72
$ws.Range("K7").Value = 219

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("K2").Value = 68
$ws.Range("K3").Value = 82
$ws.Range("K6").Value = 87
$ws.Range("K7").Value = 255

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("K2").Value = 47
$ws.Range("K3").Value = 46
$ws.Range("K6").Value = 48
$ws.Range("K7").Value = 151

$ws = $wb.Worksheets.Item("Irving Park")
$ws.Range("K6").Value = 21
$ws.Range("K7").Value = 50

$ws = $wb.Worksheets.Item("River North")
$ws.Range("K4").Value = 6
$ws.Range("K6").Value = 48
$ws.Range("K7").Value = 83

$ws = $wb.Worksheets.Item("Ashburn")
$ws.Range("K6").Value = 16
$ws.Range("K7").Value = 46

$ws = $wb.Worksheets.Item("Hermosa")
$ws.Range("K2").Value = 19
$ws.Range("K4").Value = 4
$ws.Range("K7").Value = 59

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("K3").Value = 54
$ws.Range("K6").Value = 90
$ws.Range("K7").Value = 203

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("K3").Value = 12
$ws.Range("K7").Value = 58

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("K6").Value = 12
$ws.Range("K7").Value = 53

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("K3").Value = 53
$ws.Range("K7").Value = 151

$ws = $wb.Worksheets.Item("Near South Side")
$ws.Range("K6").Value = 10
$ws.Range("K7").Value = 35

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("K4").Value = 6
$ws.Range("K7").Value = 128

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("K2").Value = 27
$ws.Range("K7").Value = 63

$ws = $wb.Worksheets.Item("Wrigleyville")
$ws.Range("K6").Value = 5
$ws.Range("K7").Value = 8

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("K4").Value = 7
$ws.Range("K6").Value = 30
$ws.Range("K7").Value = 68

$ws = $wb.Worksheets.Item("East Side")
$ws.Range("K3").Value = 10
$ws.Range("K7").Value = 25

$ws = $wb.Worksheets.Item("Kenwood")
$ws.Range("K3").Value = 14
$ws.Range("K7").Value = 39

$ws = $wb.Worksheets.Item("Wicker Park")
$ws.Range("K3").Value = 4
$ws.Range("K7").Value = 37

$ws = $wb.Worksheets.Item("Albany Park")
$ws.Range("K6").Value = 17
$ws.Range("K7").Value = 41

$ws = $wb.Worksheets.Item("West Town")
$ws.Range("K2").Value = 9
$ws.Range("K7").Value = 47

$ws = $wb.Worksheets.Item("West Elsdon")
$ws.Range("K2").Value = 6
$ws.Range("K6").Value = 15
$ws.Range("K7").Value = 25

$ws = $wb.Worksheets.Item("Pullman")
$ws.Range("K2").Value = 7
$ws.Range("K3").Value = 4
$ws.Range("K7").Value = 20

$ws = $wb.Worksheets.Item("North Park")
$ws.Range("K6").Value = 4
$ws.Range("K7").Value = 15

$ws = $wb.Worksheets.Item("Morgan Park")
$ws.Range("K2").Value = 8
$ws.Range("K7").Value = 43

$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Range("K4").Value = 6
$ws.Range("K6").Value = 25
$ws.Range("K7").Value = 58

$ws = $wb.Worksheets.Item("Clearing")
$ws.Range("K2").Value = 5
$ws.Range("K7").Value = 13

$ws = $wb.Worksheets.Item("Oakland")
$ws.Range("K3").Value = 3
$ws.Range("K7").Value = 16

$ws = $wb.Worksheets.Item("Sheffield & DePaul")
$ws.Range("K5").Value = 7
$ws.Range("K6").Value = 9

$ws = $wb.Worksheets.Item("Jackson Park")
$ws.Range("K6").Value = 4
$ws.Range("K7").Value = 6
